$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added for this market/product. It belongs at the
# top of this date-ordered block (row 113), so push the existing rows
# 113:141 down to 114:142 and populate the newly freed row 113.
$ws.Rows("113:113").Insert()

$ws.Range("A113").Value = 5
$ws.Range("B113").Value = "Macroferia Regional de Talca"
$ws.Range("C113").Value = "Maule"
$ws.Range("D113").Value = 44785
$ws.Range("E113").Value = 7
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100108
$ws.Range("H113").Value = "Tropicales y subtropicales"
$ws.Range("I113").Value = 100108002
$ws.Range("J113").Value = "Mango"
$ws.Range("K113").Value = "Sin especificar"
$ws.Range("L113").Value = "Primera"
$ws.Range("M113").Value = 228
$ws.Range("N113").Value = 10000
$ws.Range("O113").Value = 10000
$ws.Range("P113").Value = 10000
$ws.Range("Q113").Value = "$/bandeja 4 kilos"
$ws.Range("R113").Value = "Brasil"
$ws.Range("S113").Value = 2500
$ws.Range("T113").Value = 4
